# Update "具体时间范围" (column E) date ranges on every data sheet so the
# dash between the start/end timestamps is surrounded by spaces, e.g.
#   "2024.03.13 10:00-05.21 19:00"  ->  "2024.03.13 10:00 - 05.21 19:00"
# and bump a handful of "想去人数" (column F) counts to the refreshed
# numbers captured at the time the page was regenerated.

$wb = $excel.ActiveWorkbook

# Per-sheet (1-indexed sheet position) row -> new F-column value overrides.
$fChangesBySheet = @{
    1 = @{ 6 = 938; 8 = 13; 9 = 978; 10 = 773; 11 = 214; 14 = 800; 15 = 264; 16 = 567; 18 = 1313; 20 = 441; 21 = 1136; 22 = 2828; 23 = 1345; 24 = 672; 25 = 179; 26 = 1253; 28 = 988; 29 = 336; 30 = 1959; 31 = 42; 32 = 7; 33 = 1364 }
    2 = @{ 13 = 2 }
    3 = @{}
    4 = @{ 13 = 938; 16 = 13; 17 = 978; 18 = 773; 19 = 214; 27 = 800; 28 = 264; 29 = 567; 31 = 1313; 33 = 441; 34 = 1136; 35 = 2828; 36 = 1345; 37 = 672; 38 = 179; 39 = 1253; 41 = 2; 43 = 988; 44 = 336; 45 = 1959; 46 = 42; 47 = 7; 48 = 1364 }
}

for ($sheetIdx = 1; $sheetIdx -le $wb.Worksheets.Count; $sheetIdx++) {
    $ws = $wb.Worksheets.Item($sheetIdx)
    $lastRow = $ws.UsedRange.Rows.Count

    $fChanges = $fChangesBySheet[$sheetIdx]

    for ($r = 2; $r -le $lastRow; $r++) {
        # Column E = 5: "2024.MM.DD HH:MM-MM.DD HH:MM" -> add spaces around the dash.
        $eCell = $ws.Cells.Item($r, 5)
        $eVal = $eCell.Value2
        if ($eVal -ne $null -and $eVal -ne "") {
            $eCell.Value = $eVal.Replace("-", " - ")
        }

        # Column F = 6: refreshed "想去人数" counts for specific rows.
        if ($fChanges -ne $null -and $fChanges.ContainsKey($r)) {
            $ws.Cells.Item($r, 6).Value = $fChanges[$r]
        }
    }
}
